# Setting up run modes for Test data and implementing parameterization
$wb = $excel.ActiveWorkbook

$wsSuite = $wb.Worksheets.Item("test_suite")
$wsAdd   = $wb.Worksheets.Item("AddCustomerTest")
$wsOpen  = $wb.Worksheets.Item("OpenAccountTest")

# --- test_suite: flip the AddCustomerTest run flag from N to Y ---
$wsSuite.Range("B4").Value = "Y"

# --- AddCustomerTest: add a new "runmode" column (E) with Y/N flags ---
$wsAdd.Range("E1").Value = "runmode"
$wsAdd.Range("E2").Value = "Y"
$wsAdd.Range("E3").Value = "N"
$wsAdd.Range("E4").Value = "Y"
$wsAdd.Range("E5").Value = "Y"

# --- Update each sheet's selection to match where the author left off ---
$wsSuite.Range("E6").Select()
$wsOpen.Range("H22").Select()

# AddCustomerTest becomes the active tab/sheet, selection at L17
$wsAdd.Activate()
$wsAdd.Range("L17").Select()
